$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

# Append a new locator row (case_number / xpath / //td[text()="%s"]) mirroring
# the existing rows' layout (NAME, BY, VALUE columns A:C).
$ws.Range("A5").Value = "case_number"
$ws.Range("B5").Value = "xpath"
$ws.Range("C5").Value = '//td[text()="%s"]'

# Match the formatting used by the other data rows (row 4) for columns A:C.
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A5:C5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Move the active selection to the newly added cell, as in the edited file.
$ws.Range("A5").Select() | Out-Null

$wb.Save()
